# Update Data Sources from LFX (2026-02-06)
#
# Replaces the table style applied to every table in the deck from the
# custom "Medium" style ({C13E3245-A485-425E-8DFC-A9FE81DB0A83}) to the
# built-in "No Style, Table Grid" style
# ({7E67995A-638B-42D5-9D7F-B0D75F3B4303}).

$oldStyleId = "{C13E3245-A485-425E-8DFC-A9FE81DB0A83}"
$newStyleId = "{7E67995A-638B-42D5-9D7F-B0D75F3B4303}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
